# "Invoice has been Fixed" - update FullInvoice sheet test data:
#  - phone number expected value changed
#  - room count expected value changed
#  - several "ActualResult/Pass-Fail" style cells cleared (test not yet re-run)
#  - view scrolled / selection moved

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FullInvoice")

# Update expected phone number and room count
$ws.Range("B2").Value = "0123456781"
$ws.Range("E2").Value = 9

# Clear actual-result columns that no longer have cached values
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()

$ws.Range("C5").ClearContents()

$ws.Range("A8").Value = 1
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()

$ws.Range("A12").Value = 1
$ws.Range("C12").ClearContents()
$ws.Range("D12").ClearContents()

$ws.Range("A16").Value = 1
$ws.Range("E16").ClearContents()
$ws.Range("F16").ClearContents()
$ws.Range("G16").ClearContents()

# Move view / selection like in the saved file
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("H2").Select()
